$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

$ws.Range("B5").Value = "mc"
$ws.Range("D5").Value = "Is dit een goede nieuwe vraag??"
$ws.Range("E5").Value = "['A. test 1', ' B. Test 2', ' C. Test 3']"
$ws.Range("F5").Value = 1
$ws.Range("L5").Value = ""
